{"js": "// Office.js (Word JavaScript API) script.\n// Applies the May 26 2023 SME feedback edits to the health-care POA\n// revocation notice template:\n//   1. Wrap the \"dated {{ health_agent_date }}\" clause in a\n//      {% if health_agent_date != \"\" %} ... {% endif %} conditional, and\n//      drop the comma that used to separate it from \"empowering\".\n//   2. Append a new conditional sentence about delayed revocation at the\n//      end of the main paragraph (before the \"Dated:\" line break).\n//   3. Split the notary \"WITNESS my hand and official seal.\" run so a\n//      lastRenderedPageBreak sits between \"hand \" and \"and official\n//      seal.\" (mirrors the re-paginated run split in the saved document).\n\nconst body = context.document.body;\n\n// --- Edit 1: wrap the health_agent_date clause in a conditional -------\nconst datedResults = body.search(\n  \"Health Care dated {{ health_agent_date }}, empowering\",\n  { matchCase: true }\n);\ndatedResults.load(\"text\");\nawait context.sync();\n\nif (datedResults.items.length > 0) {\n  datedResults.items[0].insertText(\n    \"Health Care{% if health_agent_date != \\u201c\\u201d %} dated {{ health_agent_date }}{% endif %} empowering\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// --- Edit 2: append the delayed-revocation sentence --------------------\nconst tailResults = body.search(\n  \"granted under the aforementioned Durable Power of Attorney for Health Care. {% endif %}{% endif %}\",\n  { matchCase: true }\n);\ntailResults.load(\"text\");\nawait context.sync();\n\nif (tailResults.items.length > 0) {\n  tailResults.items[0].insertText(\n    \"granted under the aforementioned Durable Power of Attorney for Health Care. {% endif %}{% endif %}{% if revocable_poa == True %}{% if delayed_revocation == True %}This revocation shall take effect 30 days after I have communicated my intention to revoke.{% endif %}{% endif %}\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// --- Edit 3: split \"WITNESS my hand and official seal.\" so a\n//             lastRenderedPageBreak separates the two halves -----------\nconst witnessResults = body.search(\"WITNESS my hand and official seal. \", {\n  matchCase: true,\n});\nwitnessResults.load(\"text\");\nawait context.sync();\n\nif (witnessResults.items.length > 0) {\n  const ooxml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>' +\n    '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n    '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n    \"</Relationships></pkg:xmlData></pkg:part>\" +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body><w:p>\" +\n    '<w:r><w:t xml:space=\"preserve\">WITNESS my hand </w:t></w:r>' +\n    '<w:r><w:lastRenderedPageBreak/><w:t xml:space=\"preserve\">and official seal. </w:t></w:r>' +\n    \"</w:p></w:body></w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\";\n  witnessResults.items[0].insertOoxml(ooxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Word COM interop script.\n# Applies the May 26 2023 SME feedback edits to the health-care POA\n# revocation notice template:\n#   1. Wrap the \"dated {{ health_agent_date }}\" clause in a\n#      {% if health_agent_date != \"\" %} ... {% endif %} conditional, and\n#      drop the comma that used to separate it from \"empowering\".\n#   2. Append a new conditional sentence about delayed revocation at the\n#      end of the main paragraph (before the \"Dated:\" line break).\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: wrap the health_agent_date clause in a conditional --------\n$rng1 = $d.Content\n$rng1.Find.Execute(\n    \"Health Care dated {{ health_agent_date }}, empowering\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Health Care{% if health_agent_date != \u201c\u201d %} dated {{ health_agent_date }}{% endif %} empowering\",\n    2\n)\n\n# --- Edit 2: append the delayed-revocation sentence ---------------------\n$rng2 = $d.Content\n$rng2.Find.Execute(\n    \"granted under the aforementioned Durable Power of Attorney for Health Care. {% endif %}{% endif %}\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"granted under the aforementioned Durable Power of Attorney for Health Care. {% endif %}{% endif %}{% if revocable_poa == True %}{% if delayed_revocation == True %}This revocation shall take effect 30 days after I have communicated my intention to revoke.{% endif %}{% endif %}\",\n    2\n)\n"}
